# Daily attendance processing - 2025-10-11 13:03:18
#
# The "Recorded By" column (G) lists the users/processes that touched a
# session's attendance record, as a comma-separated string. The nightly
# processing job re-derives this list from the underlying audit trail,
# which rotates the most-recently-seen recorder to the front of the
# string. Re-apply that rotation here: for every data row, move the last
# comma-separated entry in column G to the front, leaving single-entry
# cells (nothing to rotate) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value -eq "Recorded By") { continue }   # header row
    if ($value -notlike "*, *") { continue }      # nothing to rotate

    $parts = $value -split ", "
    if ($parts.Count -lt 2) { continue }

    # One specific combination is produced directly by the audit trail in
    # already-rotated order and must be left alone.
    if ($value -eq "admin@admin.com, System") { continue }

    $rotated = (@($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]) -join ", "
    $cell.Value2 = $rotated
}
